# The "list view" (language table) on Sheet1 is missing its header row
# describing the tool itself ("Astronomical Processing"). Insert a new
# row right under the EN/FR/DE header row (row 3) and fill it in with the
# English / French / German captions, pushing the existing data rows
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 4, shifting rows 4:21 down to 5:22.
$ws.Rows("4:4").Insert()

# Fill in the new row - English first, then German, then French, matching
# the order the strings were originally authored in.
$ws.Range("D4").Value = "Astronomical Processing"
$ws.Range("F4").Value = "Astronomische Verarbeitung"
$ws.Range("E4").Value = "Traitement astronomique"

# Reflect the new selection/active cell left behind after the edit.
[void]$ws.Range("E22").Select()
